$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Turn the existing A1:E3 range into a real Excel Table (ListObject) ---
# Do this before writing the new formulas so the structured references
# (Tableau3[...]) resolve against the table.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:E3"), $null, 1)
$lo.Name = "Tableau3"

# --- New row labels -------------------------------------------------
# Write the "COÛT Horaire" / "COÛT Total" labels before "Mise en route"
# so the shared-string table ends up with the same insertion order as
# the authored workbook.
$ws.Range("A12").Value = "COÛT Horaire"
$ws.Range("A13").Value = "COÛT Total"
$ws.Range("A2").Value = "Mise en route"

# --- Row 2 ("Mise en route") hours, now numeric instead of "1h"/"2h" text
$ws.Range("B2").Value = 2.5
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1

# --- Row 3 ("Étude logiciels") hours, now numeric instead of text
$ws.Range("B3").Value = 2
$ws.Range("D3").Value = 0.5
$ws.Range("E3").Value = 1

# Number-format the whole hours block (B2:E3) with one decimal place,
# this also materialises the still-empty C3/E2 cells with the style.
$ws.Range("B2:E3").NumberFormat = "0.0"

# --- Row 12 ("COÛT Horaire") — hourly rate per person, in euros
$ws.Range("B12").Value = 20
$ws.Range("C12").Value = 20
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = 30
$ws.Range("B12:E12").NumberFormat = "#,##0\ ""€"";[Red]\-#,##0\ ""€"""

# --- Row 13 ("COÛT Total") — total cost per person, computed from the table
$ws.Range("B13").Formula = "=SUM(Tableau3[Nb d''h Tom])*B12"
$ws.Range("C13").Formula = "=SUM(Tableau3[Nb d''h Loïc])*C12"
$ws.Range("D13").Formula = "=SUM(Tableau3[Nb d''h Anuar])*D12"
$ws.Range("E13").Formula = "=SUM(Tableau3[Nb d''h M. KARINE])*E12"
$ws.Range("B13:E13").NumberFormat = "#,##0\ ""€"""

# --- Cosmetic sheet tweaks -------------------------------------------
# Column C widened to match column B instead of auto best-fit.
$ws.Columns.Item(3).ColumnWidth = 9.5

# Page setup: A4 portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# View: zoomed in to 130%, selection moved to E6.
$ws.Application.ActiveWindow.Zoom = 130
$null = $ws.Range("E6").Select()
